# Update column F ('想去人数' / interest counts) across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F12").Value = 1503
$ws.Range("F14").Value = 1460
$ws.Range("F17").Value = 306
$ws.Range("F20").Value = 1034
$ws.Range("F24").Value = 1464
$ws.Range("F26").Value = 144
$ws.Range("F29").Value = 1100
$ws.Range("F31").Value = 991
$ws.Range("F35").Value = 1071
$ws.Range("F42").Value = 1619
$ws.Range("F44").Value = 57

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 1461
$ws.Range("F11").Value = 1194
$ws.Range("F15").Value = 29
$ws.Range("F19").Value = 441
$ws.Range("F20").Value = 22
$ws.Range("F32").Value = 29
$ws.Range("F39").Value = 49
$ws.Range("F40").Value = 49

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2824
$ws.Range("F6").Value = 4553
$ws.Range("F10").Value = 674
$ws.Range("F12").Value = 260
$ws.Range("F13").Value = 879
$ws.Range("F14").Value = 227
$ws.Range("F15").Value = 519

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2824
$ws.Range("F5").Value = 4553
$ws.Range("F6").Value = 674
$ws.Range("F8").Value = 260
$ws.Range("F9").Value = 260
$ws.Range("F10").Value = 879
$ws.Range("F11").Value = 879
$ws.Range("F14").Value = 1461
$ws.Range("F18").Value = 1194
$ws.Range("F19").Value = 1503
$ws.Range("F21").Value = 1460
$ws.Range("F27").Value = 1034
$ws.Range("F29").Value = 519
$ws.Range("F30").Value = 519
$ws.Range("F31").Value = 441
$ws.Range("F32").Value = 1464
$ws.Range("F33").Value = 144
$ws.Range("F36").Value = 1100
$ws.Range("F38").Value = 991
$ws.Range("F40").Value = 1071
$ws.Range("F47").Value = 1619
$ws.Range("F50").Value = 49

